$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Requisitos" (Requirements) list had two entries:
#   row 23 (B23/C23): LOT2028 - Tecnologia de Processos Fermentativos (Requisito fraco)
#   row 24 (B24/C24): LOT2052 - Tecnologia de Bebidas Experimental (Indicacao de Conjunto)
# The edit reorders the shared-string table so the LOT2052 entry now comes
# first, i.e. the text content of these two rows is swapped.

$lot2028 = "LOT2028 -  Tecnologia de Processos Fermentativos  (Requisito fraco)`n"
$lot2052 = "LOT2052 -  Tecnologia de Bebidas Experimental  (Indicação de Conjunto)`n"

$ws.Range("B23").Value = $lot2052
$ws.Range("C23").Value = $lot2052

$ws.Range("B24").Value = $lot2028
$ws.Range("C24").Value = $lot2028
